$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. New label / header text, entered in the order that reproduces
#    the target shared-string table ordering (indices 3-12).
# -----------------------------------------------------------------
$ws.Range("K17").Value2 = "accuracy"
$ws.Range("K18").Value2 = "time"

$ws.Range("G12").Value2 = "mean"
$ws.Range("G13").Value2 = "SD"

$ws.Range("G1").Value2  = "(1.4.3) 90/10 split, 10 runs"
$ws.Range("K1").Value2  = "(1.4.4)90/10, gaussian sig=0.55"
$ws.Range("O1").Value2  = "(1.4.4)90/10, gaussian sig=0.75"

$ws.Range("K15").Value2 = "many persons, 90/10 split"

$ws.Range("D1").Value2  = "50/50 split"
$ws.Range("D2").Value2  = "k testing"

# Re-use ("mean" / "SD") for the other two mini result tables
$ws.Range("K12").Value2 = "mean"
$ws.Range("O12").Value2 = "mean"
$ws.Range("L16").Value2 = "mean"

$ws.Range("K13").Value2 = "SD"
$ws.Range("O13").Value2 = "SD"
$ws.Range("M16").Value2 = "SD"

# -----------------------------------------------------------------
# 2. New formulas for the third (O/P/Q) mean/SD block
# -----------------------------------------------------------------
$ws.Range("P12").Formula = "=AVERAGE(P2:P11)"
$ws.Range("Q12").Formula = "=AVERAGE(Q2:Q11)"
$ws.Range("P13").Formula = "=STDEV(P2:P12)"
$ws.Range("Q13").Formula = "=STDEV(Q2:Q12)"

# -----------------------------------------------------------------
# 3. New small "overall accuracy/time" summary table (K17:M18)
# -----------------------------------------------------------------
$ws.Range("A17").Copy()
$ws.Range("K17").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K17").Value2 = "accuracy"

$ws.Range("A18").Copy()
$ws.Range("K18").PasteSpecial(-4122)
$ws.Range("K18").Value2 = "time"

$ws.Range("B17").Copy()
$ws.Range("L17").PasteSpecial(-4122)
$ws.Range("L17").Value2 = 93.15

$ws.Range("B17").Copy()
$ws.Range("M17").PasteSpecial(-4122)
$ws.Range("M17").Value2 = 0.72140000000000004
$ws.Range("M17").Borders.LineStyle = -4142   # xlLineStyleNone

$ws.Range("B18").Copy()
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("L18").Value2 = 8563.25

$ws.Range("B18").Copy()
$ws.Range("M18").PasteSpecial(-4122)
$ws.Range("M18").Value2 = 42.113
$ws.Range("M18").Borders.LineStyle = -4142

$ws.Range("A1").Copy()   # clear clipboard marquee

# -----------------------------------------------------------------
# 4. Centre-align + merge the three grouped headers in row 1 and the
#    new grouped header in row 15
# -----------------------------------------------------------------
$ws.Range("G1:I1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("K1:M1").HorizontalAlignment = -4108
$ws.Range("O1:Q1").HorizontalAlignment = -4108
$ws.Range("K15:M15").HorizontalAlignment = -4108

$ws.Range("G1:I1").Merge()
$ws.Range("K1:M1").Merge()
$ws.Range("O1:Q1").Merge()
$ws.Range("K15:M15").Merge()

# -----------------------------------------------------------------
# 5. Column widths for the new D / E columns
# -----------------------------------------------------------------
$ws.Range("D1:E1").ColumnWidth = 12.65

# -----------------------------------------------------------------
# 6. View: scroll back to top-left and move the selection
# -----------------------------------------------------------------
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("D7").Select()
